# Add a new "Avvocato/Notaio" lookup row to the tab_Tipo_tribunale sheet.
# The new entry (ID 4) is inserted as row 5, pushing the former row 5
# ("Non conosciuto", ID 9) down to the new row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Results")

# Populate the new row 6 with the values that used to live in row 5.
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = "Non conosciuto"

# Match row 6's formatting (including the still-blank D/E/F cells) to row 5.
$ws.Range("A5:F5").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Overwrite row 5 in place with the new "Avvocato/Notaio" entry (ID 4).
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Avvocato/Notaio"

# Extend the used range / selection to reflect the newly added row.
$ws.Range("C6").Select()
